$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # G2: 22.33 -> "不可售"
    $ws.Range("G2").Value = "不可售"

    # G3: 55 -> 65
    $ws.Range("G3").Value = 65

    # F6: 129 -> 131
    $ws.Range("F6").Value = 131

    if ($name -eq "展览") {
        # F8: 4788 -> 4810
        $ws.Range("F8").Value = 4810
        # F10: 5143 -> 5157
        $ws.Range("F10").Value = 5157
        # F12: 1286 -> 1287
        $ws.Range("F12").Value = 1287
    } else {
        # 全部类型 sheet has an extra row, shifting indices by 1
        # F9: 4788 -> 4810
        $ws.Range("F9").Value = 4810
        # F11: 5143 -> 5157
        $ws.Range("F11").Value = 5157
        # F13: 1286 -> 1287
        $ws.Range("F13").Value = 1287
    }
}
